# Config is now done over the DHCP management port:
# - drop the "ansible_host_subnet" and "ansible_host_dg" columns (D:E) from
#   Sheet1 entirely, since the subnet/default-gateway are no longer
#   installed manually
# - renumber the pt-sw-accessR3/R4/R5 management IPs onto the new dhcp range

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D:E").Delete()

$ws.Range("C8").Value = "10.9.12.26"
$ws.Range("C9").Value = "10.9.12.17"
$ws.Range("C10").Value = "10.9.12.80"
